# imf_ppp/2008.xlsx data cleaning pass:
#  - drop the "Swedish" and "Uzbek" rows
#  - sort the remaining languages in descending order by their 2008 value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Swedish" row (row 18) and the "Uzbek" row (originally row 22,
# which becomes row 21 once the Swedish row above it has been removed).
# Deleting entire rows shifts everything below up and shrinks the sheet's
# used range / dimension accordingly, while leaving cell formatting intact.
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(21).Delete()

# Final ordering of the remaining 20 languages, sorted descending by value.
$data = @(
    @("English", 23.73011951594644),
    @("Chinese", 13.19492681339658),
    @("Spanish", 7.418505234518354),
    @("Japanese", 5.373811535534166),
    @("Arabic", 5.206363688915678),
    @("German", 4.738736492889176),
    @("Russian", 3.827786504167395),
    @("Portuguese", 3.560677449720673),
    @("French", 2.941539099064483),
    @("Malay-Indonesian", 2.793783400778675),
    @("Italian", 2.553072751354447),
    @("Korean", 1.659717354092937),
    @("Persian", 1.547009816542462),
    @("Turkish", 1.44433974040968),
    @("Dutch", 1.405226822571314),
    @("Thai", 0.966340895938355),
    @("Polish", 0.8841498482305022),
    @("Urdu", 0.8804951220189124),
    @("Vietnamese", 0.4888291112831661),
    @("Bengali", 0.4833908691998375)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}
